$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the task row "Crear una interfaz que permita al usuario indicar si quiere
# guardar el ingrediente buscado" (old row 20) - rows below shift up automatically.
$ws.Rows("20:20").Delete(-4162)

# --- BÚSQUEDA DE INGREDIENTES header row (row 12): add ESTADO / TIEMPO ESTIMADO ---
$ws.Range("C7").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = "HECHO"

$ws.Range("D12").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "1 hora"

# --- Fix formatting of C14 (empty "HECHO" cell, fill changes from plain to green) ---
$ws.Range("C19").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- SELECCIÓN DE INGREDIENTES header row (row 18): add ESTADO / TIEMPO ESTIMADO / TIEMPO INVERTIDO ---
$ws.Range("C7").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "HECHO"

$ws.Range("D18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "1 hora"

$ws.Range("D18").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value = "0,5 horas"

# --- back to row 12: TIEMPO INVERTIDO ---
$ws.Range("D12").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = "2 horas"

# --- Fix formatting of C20 (previously C21 before the row deletion above) ---
$ws.Range("C19").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("F18").Select() | Out-Null
